# "next step to do." - update the work-plan sheet:
#  - C3 becomes the short status note "js编写完成" with a new
#    left/top-aligned style.
#  - C2 is narrowed to just the JS-writing task line.
#  - D2 picks up the newly split-out list of remaining tasks
#    (info-list classes, jsp rewrite, wiring up interaction).
#
# Shared strings are appended in the order the cells are written, so the
# order of the three assignments below matters (it reproduces the
# uniqueCount 6 -> 8 growth with index 5/6/7 landing on the same cells
# as in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C3: new short "done" note, left/top aligned (new cell style).
$c3 = $ws.Cells.Item(3, 3)
$c3.Value = "js编写完成"
$c3.HorizontalAlignment = -4131   # xlLeft
$c3.VerticalAlignment = -4160     # xlTop

# C2: keep only the "writing the JS" line (trailing newline preserved).
$ws.Cells.Item(2, 3).Value = "主页、评论页面、购物车js的编写`n"

# D2: the newly-added list of remaining work items.
$ws.Cells.Item(2, 4).Value = "评论信息列表类`n商品信息类`n购物车商品信息类`n主页、评论页面、购物车用jsp写`n实现交互"
